# Automatische test-sync: 2025-06-26 23:08:50
# Adds the new "Wanneer zijn jullie open?" test-mail row (row 32) to the
# "Logs" sheet, extends the conditional formatting ranges that covered
# rows 2-31 so they also cover the new row 32, and bumps the
# "Openingstijden / Locatie" tally on the "Dashboard" sheet from 8 to 9.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log row (row 32) ------------------------------------
$logs.Range("A32").Value = "Wanneer zijn jullie open?"
$logs.Range("B32").Value = "mailmind.test@zohomail.eu"
$logs.Range("C32").Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D32").Value = "Openingstijden / Locatie"
$logs.Range("E32").Value = "Beste klant,`n`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`n`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F32").Value = "2025-06-26 23:08:37"
$logs.Range("G32").Value = "Ja"
$logs.Range("H32").Value = "Nee"
$logs.Range("I32").Value = "Ja"

# Writing the multi-line E32 text makes the host auto-mark the row with a
# custom height; AutoFit() settles it back to the sheet's implicit default
# (no row before it in the sheet carries an explicit height either).
$logs.Rows.Item(32).AutoFit()

# --- 2. Extend the conditional formatting ranges to include row 32 ---------
# The workbook has four conditional-formatting blocks (columns D, G, H, I),
# each spanning rows 2-31. Recreate every rule on the extended row 2-32
# range so the formatting that highlights "Klacht"/"Ja"/... keeps applying
# to the freshly appended row.

$dRange = $logs.Range("D2:D32")
$dRange.FormatConditions.Delete()
$d1 = $dRange.FormatConditions.Add(1, 3, '="Klacht"');             $d1.Priority = 1
$d2 = $dRange.FormatConditions.Add(1, 3, '="Bestelling"');         $d2.Priority = 2
$d3 = $dRange.FormatConditions.Add(1, 3, '="Informatieaanvraag"'); $d3.Priority = 3
$d4 = $dRange.FormatConditions.Add(1, 3, '="Afmelding"');          $d4.Priority = 4
$d5 = $dRange.FormatConditions.Add(1, 3, '="Retour"');             $d5.Priority = 5
$d6 = $dRange.FormatConditions.Add(1, 3, '="Overig"');             $d6.Priority = 6

$gRange = $logs.Range("G2:G32")
$gRange.FormatConditions.Delete()
$g1 = $gRange.FormatConditions.Add(1, 3, '="Ja"');  $g1.Priority = 7
$g2 = $gRange.FormatConditions.Add(1, 3, '="Nee"'); $g2.Priority = 8

$hRange = $logs.Range("H2:H32")
$hRange.FormatConditions.Delete()
$h1 = $hRange.FormatConditions.Add(1, 3, '="Ja"'); $h1.Priority = 9

$iRange = $logs.Range("I2:I32")
$iRange.FormatConditions.Delete()
$i1 = $iRange.FormatConditions.Add(1, 3, '="Ja"'); $i1.Priority = 10

# --- 3. Update the Dashboard tally for "Openingstijden / Locatie" ----------
$dashboard.Range("B3").Value = 9
